# Apply the "Updated cryptos list" price/volume refresh (and the
# RenderToken/TheSandbox rank swap at rows 37-38) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.226.90"
$ws.Cells.Item(2, 5).Value = "  +0.64%  "
$ws.Cells.Item(3, 4).Value = "1.905.05"
$ws.Cells.Item(3, 5).Value = "  +0.66%  "
$ws.Cells.Item(4, 4).Value = "'1.001"
$ws.Cells.Item(4, 5).Value = "  -0.21%  "
$ws.Cells.Item(5, 4).Value = "'306.19"
$ws.Cells.Item(5, 5).Value = "  -0.16%  "
$ws.Cells.Item(6, 4).Value = "'1.000"
$ws.Cells.Item(6, 5).Value = "  -0.17%  "
$ws.Cells.Item(7, 4).Value = "'0.5419"
$ws.Cells.Item(7, 5).Value = "  +3.68%  "
$ws.Cells.Item(8, 4).Value = "'0.3807"
$ws.Cells.Item(8, 5).Value = "  +1.34%  "
$ws.Cells.Item(9, 4).Value = "'0.07289"
$ws.Cells.Item(9, 5).Value = "  +0.37%  "
$ws.Cells.Item(10, 4).Value = "'22.09"
$ws.Cells.Item(10, 5).Value = "  +4.75%  "
$ws.Cells.Item(11, 4).Value = "'0.9023"
$ws.Cells.Item(11, 5).Value = "  +0.38%  "
$ws.Cells.Item(12, 4).Value = "'0.08201"
$ws.Cells.Item(12, 5).Value = "  +0.40%  "
$ws.Cells.Item(13, 4).Value = "'95.67"
$ws.Cells.Item(13, 5).Value = "  -0.54%  "
$ws.Cells.Item(14, 4).Value = "'5.349"
$ws.Cells.Item(14, 5).Value = "  +1.08%  "
$ws.Cells.Item(15, 4).Value = "'1.001"
$ws.Cells.Item(15, 5).Value = "  -0.16%  "
$ws.Cells.Item(16, 4).Value = "'14.83"
$ws.Cells.Item(16, 5).Value = "  +1.81%  "
$ws.Cells.Item(17, 4).Value = "'0.000008640"
$ws.Cells.Item(17, 5).Value = "  +0.72%  "
$ws.Cells.Item(18, 4).Value = "'1.0000"
$ws.Cells.Item(18, 5).Value = "  -0.29%  "
$ws.Cells.Item(19, 4).Value = "1.241.57"
$ws.Cells.Item(19, 5).Value = "  -34.63%  "
$ws.Cells.Item(20, 4).Value = "27.262.51"
$ws.Cells.Item(20, 5).Value = "  +0.60%  "
$ws.Cells.Item(21, 4).Value = "'5.048"
$ws.Cells.Item(21, 5).Value = "  -0.55%  "
$ws.Cells.Item(22, 4).Value = "'10.81"
$ws.Cells.Item(22, 5).Value = "  +1.20%  "
$ws.Cells.Item(23, 4).Value = "'6.513"
$ws.Cells.Item(23, 5).Value = "  +1.70%  "
$ws.Cells.Item(24, 4).Value = "'148.31"
$ws.Cells.Item(24, 5).Value = "  -0.23%  "
$ws.Cells.Item(25, 4).Value = "'2.299"
$ws.Cells.Item(25, 5).Value = "  +0.43%  "
$ws.Cells.Item(26, 4).Value = "'18.36"
$ws.Cells.Item(26, 5).Value = "  +1.04%  "
$ws.Cells.Item(27, 4).Value = "'1.757"
$ws.Cells.Item(27, 5).Value = "  +1.39%  "
$ws.Cells.Item(28, 4).Value = "'116.81"
$ws.Cells.Item(28, 5).Value = "  +1.62%  "
$ws.Cells.Item(29, 4).Value = "'4.851"
$ws.Cells.Item(29, 5).Value = "  +1.50%  "
$ws.Cells.Item(30, 4).Value = "'4.660"
$ws.Cells.Item(30, 5).Value = "  -3.74%  "
$ws.Cells.Item(31, 4).Value = "'0.09206"
$ws.Cells.Item(31, 5).Value = "  -0.16%  "
$ws.Cells.Item(32, 4).Value = "'0.8258"
$ws.Cells.Item(32, 5).Value = "  +4.65%  "
$ws.Cells.Item(33, 4).Value = "'0.05065"
$ws.Cells.Item(33, 5).Value = "  +0.58%  "
$ws.Cells.Item(34, 4).Value = "'1.220"
$ws.Cells.Item(34, 5).Value = "  +0.86%  "
$ws.Cells.Item(35, 4).Value = "'3.012"
$ws.Cells.Item(35, 5).Value = "  +1.12%  "
$ws.Cells.Item(36, 4).Value = "'3.323"
$ws.Cells.Item(36, 5).Value = "  -3.11%  "
$ws.Cells.Item(37, 2).Value = "RenderToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(37, 4).Value = "'2.695"
$ws.Cells.Item(37, 5).Value = "  +3.56%  "
$ws.Cells.Item(38, 2).Value = "TheSandbox"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(38, 4).Value = "'0.6023"
$ws.Cells.Item(38, 5).Value = "  +5.26%  "
$ws.Cells.Item(39, 4).Value = "'0.01995"
$ws.Cells.Item(39, 5).Value = "  +0.53%  "
$ws.Cells.Item(40, 5).Value = "  +0.08%  "
$ws.Cells.Item(41, 4).Value = "'9.254"
$ws.Cells.Item(41, 5).Value = "  +2.35%  "
$ws.Cells.Item(42, 4).Value = "'6.657"
$ws.Cells.Item(42, 5).Value = "  +1.62%  "
$ws.Cells.Item(43, 4).Value = "'115.92"
$ws.Cells.Item(43, 5).Value = "  -0.33%  "
$ws.Cells.Item(44, 4).Value = "'0.5140"
$ws.Cells.Item(44, 5).Value = "  +5.64%  "
$ws.Cells.Item(45, 4).Value = "'0.1530"
$ws.Cells.Item(45, 5).Value = "  +1.03%  "
$ws.Cells.Item(46, 4).Value = "'10.21"
$ws.Cells.Item(46, 5).Value = "  +1.65%  "
$ws.Cells.Item(47, 4).Value = "'0.9998"
$ws.Cells.Item(47, 5).Value = "  -0.18%  "
$ws.Cells.Item(48, 4).Value = "'1.639"
$ws.Cells.Item(48, 5).Value = "  +0.99%  "
$ws.Cells.Item(49, 4).Value = "'38.09"
$ws.Cells.Item(49, 5).Value = "  -0.13%  "
$ws.Cells.Item(50, 5).Value = "  +2.89%  "
$ws.Cells.Item(51, 4).Value = "'63.54"
$ws.Cells.Item(51, 5).Value = "  +0.06%  "
